# Update metric names in column A to new naming convention
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "selfemployedRate"
$ws.Range("A4").Value = "unemployedRate"
$ws.Range("A7").Value = "selfemployed"
$ws.Range("A8").Value = "unemployed"
$ws.Range("A9").Value = "inactive"
$ws.Range("A18").Value = "L3PlusPerc"
$ws.Range("A21").Value = "employmentProjection"

# Move selection/view to bottom of sheet (row 23/24) as in target workbook
$ws.Activate()
$ws.Range("A24").Select()
$excel.ActiveWindow.ScrollRow = 23
